# Update league data: cyclically rotate the row data (columns B:AD) among
# rows 172, 173 and 174 (row index/id column A is left untouched):
#   new row172 <- old row174
#   new row173 <- old row172
#   new row174 <- old row173

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colFirst = 2   # column B
$colLast  = 30  # column AD

# Capture current ("before") values for the three rows across columns B:AD.
$row172 = @{}
$row173 = @{}
$row174 = @{}
for ($c = $colFirst; $c -le $colLast; $c++) {
    $row172[$c] = $ws.Cells.Item(172, $c).Value2
    $row173[$c] = $ws.Cells.Item(173, $c).Value2
    $row174[$c] = $ws.Cells.Item(174, $c).Value2
}

# Write back the rotated values.
for ($c = $colFirst; $c -le $colLast; $c++) {
    $ws.Cells.Item(172, $c).Value2 = $row174[$c]
    $ws.Cells.Item(173, $c).Value2 = $row172[$c]
    $ws.Cells.Item(174, $c).Value2 = $row173[$c]
}
